# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    and fill it with the quarterly fund-holding detail.
# 2. Prepend a "2022-Q1" summary row to the "总计" (grand total) sheet and
#    renumber the existing rows.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force a value that *looks* numeric (e.g. "001304" or "2.13") to be
    # stored as literal text instead of being auto-coerced to a number,
    # then drop the temporary "@" number-format so the cell is left with
    # no explicit style, matching a freshly authored text cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet
# ---------------------------------------------------------------------

$q4sheet = $wb.Worksheets.Item("2021-Q4")

# Use the "2021-Q4" sheet (same column layout: 基金代码/基金名称/基金规模/
# 股票总仓位/仓位占比/持有市值(亿元)/仓位排名) as a formatting template so the
# new sheet picks up identical header/row styles without inventing new ones.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4sheet)
$newSheet.Name = "2022-Q1"

$q4sheet.Range("A1:H2").Copy($newSheet.Range("A1"))
$newSheet.Range("A1").ClearContents()

# Stamp out the template row for the remaining four data rows (3-6).
$q4sheet.Range("A2:H2").Copy($newSheet.Range("A3:H3"))
$q4sheet.Range("A2:H2").Copy($newSheet.Range("A4:H4"))
$q4sheet.Range("A2:H2").Copy($newSheet.Range("A5:H5"))
$q4sheet.Range("A2:H2").Copy($newSheet.Range("A6:H6"))

# Header row (plain Value assignment keeps the bold/bordered "s=2" style
# that was already copied in from the template row).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 001304 建信鑫安回报灵活配置混合
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "001304"
Set-TextValue $newSheet.Range("C2") "建信鑫安回报灵活配置混合"
Set-TextValue $newSheet.Range("D2") "2.13"
Set-TextValue $newSheet.Range("E2") "66.83"
Set-TextValue $newSheet.Range("F2") "6.87"
Set-TextValue $newSheet.Range("G2") "0.1463"
$newSheet.Range("H2").Value = 1

# Row 3 - 006279 中金瑞祥灵活配置混合A
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "006279"
Set-TextValue $newSheet.Range("C3") "中金瑞祥灵活配置混合A"
Set-TextValue $newSheet.Range("D3") "2.10"
Set-TextValue $newSheet.Range("E3") "59.54"
Set-TextValue $newSheet.Range("F3") "6.51"
Set-TextValue $newSheet.Range("G3") "0.1367"
$newSheet.Range("H3").Value = 4

# Row 4 - 002585 建信兴利灵活配置混合
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "002585"
Set-TextValue $newSheet.Range("C4") "建信兴利灵活配置混合"
Set-TextValue $newSheet.Range("D4") "2.04"
Set-TextValue $newSheet.Range("E4") "61.22"
Set-TextValue $newSheet.Range("F4") "6.34"
Set-TextValue $newSheet.Range("G4") "0.1293"
$newSheet.Range("H4").Value = 2

# Row 5 - 005396 中金丰硕混合
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "005396"
Set-TextValue $newSheet.Range("C5") "中金丰硕混合"
Set-TextValue $newSheet.Range("D5") "1.83"
Set-TextValue $newSheet.Range("E5") "71.47"
Set-TextValue $newSheet.Range("F5") "6.35"
Set-TextValue $newSheet.Range("G5") "0.1162"
$newSheet.Range("H5").Value = 6

# Row 6 - 006280 中金瑞祥灵活配置混合C
$newSheet.Range("A6").Value = 4
Set-TextValue $newSheet.Range("B6") "006280"
Set-TextValue $newSheet.Range("C6") "中金瑞祥灵活配置混合C"
Set-TextValue $newSheet.Range("D6") "0.00"
Set-TextValue $newSheet.Range("E6") "59.54"
Set-TextValue $newSheet.Range("F6") "6.51"
$newSheet.Range("G6").Value = 0
$newSheet.Range("H6").Value = 4

# ---------------------------------------------------------------------
# Step 2: update the "总计" (grand total) sheet
# ---------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Shift the existing three data rows (2-4) down to (3-5) by copying
# bottom-up, which keeps styling intact without leaving stray row formats
# behind the way Rows.Insert() does.
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# New row 2 - 2022-Q1 summary
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.53

# Renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
